# Add Skill1/Skill2/Skill3 columns (one per hero position) to the Player
# "Property1" sheet, inserted right after HeroPos3 (old column N) and before
# VIPLevel (old column O), shifting all subsequent columns right by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert 3 new columns before the old column O (VIPLevel) ---------------
$ws.Columns("O:Q").Insert()

# Match column width of the new columns to column N (the HeroPos1-3 block
# they conceptually belong to).
$ws.Columns("O:Q").ColumnWidth = $ws.Columns("N").ColumnWidth()

# --- Row 1 (header/titles) ---------------------------------------------
$ws.Range("O1").Value = "Skill1"
$ws.Range("P1").Value = "Skill2"
$ws.Range("Q1").Value = "Skill3"

# --- Row 2 (Type row) ----------------------------------------------------
$ws.Range("O2").Value = "string"
$ws.Range("P2").Value = "string"
$ws.Range("Q2").Value = "string"

# --- Row 3 (Public) -------------------------------------------------------
$ws.Range("O3").Value = $ws.Range("N3").Value()
$ws.Range("P3").Value = $ws.Range("N3").Value()
$ws.Range("Q3").Value = $ws.Range("N3").Value()

# --- Row 4 (Private) -------------------------------------------------------
$ws.Range("O4").Value = $ws.Range("N4").Value()
$ws.Range("P4").Value = $ws.Range("N4").Value()
$ws.Range("Q4").Value = $ws.Range("N4").Value()

# --- Row 5 (Save) -------------------------------------------------------
$ws.Range("O5").Value = $ws.Range("N5").Value()
$ws.Range("P5").Value = $ws.Range("N5").Value()
$ws.Range("Q5").Value = $ws.Range("N5").Value()

# --- Row 6 (Cache) -------------------------------------------------------
$ws.Range("O6").Value = $ws.Range("N6").Value()
$ws.Range("P6").Value = $ws.Range("N6").Value()
$ws.Range("Q6").Value = $ws.Range("N6").Value()

# --- Row 7 (Ref) -------------------------------------------------------
$ws.Range("O7").Value = $ws.Range("N7").Value()
$ws.Range("P7").Value = $ws.Range("N7").Value()
$ws.Range("Q7").Value = $ws.Range("N7").Value()

# --- Row 8 (Upload) -------------------------------------------------------
$ws.Range("O8").Value = $ws.Range("N8").Value()
$ws.Range("P8").Value = $ws.Range("N8").Value()
$ws.Range("Q8").Value = $ws.Range("N8").Value()

# --- Row 9 (Desc) -------------------------------------------------------
$ws.Range("O9").Value = "Hero1"
$ws.Range("P9").Value = "Hero2"
$ws.Range("Q9").Value = "Hero3"

# --- Data validation sqref fix-up (K6:N6 -> K6:Q6) ------------------------
$dv = $ws.Range("K6:Q6")
$dv.Validation.Delete()
$dv.Validation.Add(3, 1, 1, "TRUE,FALSE")
$ws.Range("K6:N6").Validation.Delete()

# --- Sheet view restoration ------------------------------------------------
$sv = $ws.Application.ActiveWindow
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A12").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("Q2").Select()
